# Weekly fruit/hortaliza price update:
# Insert two new rows at the top of this sub-range (rows 43-44) with the
# latest week's data, pushing the previous rows (old 43-47) down to 45-49.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two rows before row 43, shifting existing rows 43-47 down to 45-49.
$ws.Rows.Item(43).Insert()
$ws.Rows.Item(43).Insert()

# New row 43: Sin especificar / Primera
$ws.Range("A43").Value = 3
$ws.Range("B43").Value = "Femacal de La Calera"
$ws.Range("C43").Value = "Coquimbo"
$ws.Range("D43").Value = 45173
$ws.Range("E43").Value = 5
$ws.Range("F43").Value = 100112043
$ws.Range("G43").Value = "Pepino dulce"
$ws.Range("H43").Value = "Sin especificar"
$ws.Range("I43").Value = "Primera"
$ws.Range("J43").Value = 78
$ws.Range("K43").Value = 25000
$ws.Range("L43").Value = 25000
$ws.Range("M43").Value = 25000
$ws.Range("N43").Value = "$/caja 15 kilos"
$ws.Range("O43").Value = "Provincia de Limarí"
$ws.Range("P43").Value = 1667
$ws.Range("Q43").Value = 15
$ws.Range("R43").Value = "Hortaliza"

# New row 44: Sin especificar / Segunda
$ws.Range("A44").Value = 3
$ws.Range("B44").Value = "Femacal de La Calera"
$ws.Range("C44").Value = "Coquimbo"
$ws.Range("D44").Value = 45173
$ws.Range("E44").Value = 5
$ws.Range("F44").Value = 100112043
$ws.Range("G44").Value = "Pepino dulce"
$ws.Range("H44").Value = "Sin especificar"
$ws.Range("I44").Value = "Segunda"
$ws.Range("J44").Value = 75
$ws.Range("K44").Value = 18000
$ws.Range("L44").Value = 18000
$ws.Range("M44").Value = 18000
$ws.Range("N44").Value = "$/caja 15 kilos"
$ws.Range("O44").Value = "Provincia de Limarí"
$ws.Range("P44").Value = 1200
$ws.Range("Q44").Value = 15
$ws.Range("R44").Value = "Hortaliza"
